$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (the old A3 -> www-q hyperlink) before changing values
$ws.Hyperlinks.Delete()

# Update cell values
# A2: devUrl now points to the "dev2" host instead of duplicating the prod URL
$ws.Range("A2").Value = "https://dev2.abbviepro.com/de/de.html"
# B2: prodUrl stays as the abbviepro prod URL
$ws.Range("B2").Value = "https://www.abbviepro.com/de/de.html"
# A3: now mirrors the prod URL (was previously the www-q url)
$ws.Range("A3").Value = "https://www.abbviepro.com/de/de.html"
# B3: unchanged value, stays the prod URL
$ws.Range("B3").Value = "https://www.abbviepro.com/de/de.html"

# Re-create hyperlinks for each URL cell, each pointing to its own cell text
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.abbviepro.com/de/de.html")
$ws.Hyperlinks.Add($ws.Range("A2"), "https://dev2.abbviepro.com/de/de.html")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.abbviepro.com/de/de.html")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.abbviepro.com/de/de.html")

# Apply the Hyperlink cell style to all four link cells (Add() sets its own
# style variant; reset to the shared built-in "Hyperlink" style)
$ws.Range("A2:B3").Style = "Hyperlink"

# Update the selection to match the new active cell / selected range
$ws.Range("B3").Select()
